$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cities")

# Remove the rows for cities that are not part of the reduced ICA1 dataset.
# (Deleted from the bottom up so earlier row numbers stay valid.)
$ws.Rows(26).Delete()   # Frankfurt
$ws.Rows(23).Delete()   # Magdeburg
$ws.Rows(18).Delete()   # Salzburg
$ws.Rows(16).Delete()   # Verona
$ws.Rows(13).Delete()   # Graz
$ws.Rows(10).Delete()   # Tirane
$ws.Rows(9).Delete()    # Belgrade
$ws.Rows(7).Delete()    # Bratislava
$ws.Rows(6).Delete()    # Dresden
$ws.Rows(3).Delete()    # Gdansk
$ws.Rows(2).Delete()    # Poznan

# Brno (now row 11) becomes the "solve it yourself" row: move it to the
# bottom of the table and blank out its id, turning it into the task row.
$ws.Range("A11:B11").ClearContents()
$ws.Rows(11).Delete()
$ws.Range("B15").Value = "Brno"

# Make Cities the active sheet/tab, with C15 selected.
$ws.Activate()
$ws.Range("C15").Select()
